# "updated on 20 Nov and also started rice tracking"
# Fill in day-6 ("L" column) inventory-level readings for the tongue-scraper
# tracking sheet. The downstream R ("day5-day6" sold) / S ("day6-bleed")
# formula columns, and the R18/S18 daily-average formulas, recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$daySixValues = @{
    "L8"  = 74
    "L9"  = 266
    "L10" = 1050
    "L11" = 294
    "L12" = 0
    "L13" = 69
    "L14" = 61
}

foreach ($addr in $daySixValues.Keys) {
    $ws.Range($addr).Value = $daySixValues[$addr]
}

# Preserve the scroll position / active cell reported in the author's
# session (topLeftCell moved from G2 to I2, active cell from K14 to L14).
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("L14").Select()
